# Update cryptos list: refresh prices and volume(1h) percentages;
# also re-sort rows 50-51 (Mantle now above Filecoin).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '76.217.25'
$ws.Range('E2').Value = '  +1.81%  '
$ws.Range('D3').Value = '2.918.49'
$ws.Range('E3').Value = '  +3.43%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '200.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '599.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -0.90%  '
$ws.Range('E9').Value = '  +4.39%  '
$ws.Range('D10').Value = '2.918.09'
$ws.Range('E10').Value = '  +3.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.430'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +16.75%  '
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.88'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.21%  '
$ws.Range('D14').Value = '3.455.69'
$ws.Range('E14').Value = '  +3.44%  '
$ws.Range('D15').Value = '76.044.05'
$ws.Range('E15').Value = '  +1.64%  '
$ws.Range('E16').Value = '  +2.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.60'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.43%  '
$ws.Range('D18').Value = '2.918.50'
$ws.Range('E18').Value = '  +3.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.02'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.76'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '372.01'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.51%  '
$ws.Range('E22').Value = '  +2.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +6.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.34'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').Value = '3.061.73'
$ws.Range('E26').Value = '  +3.27%  '
$ws.Range('E27').Value = '  +1.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.72'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.74%  '
$ws.Range('E29').Value = '  +6.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('E31').Value = '  +1.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '503.38'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.73'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.35%  '
$ws.Range('E34').Value = '  +1.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '165.55'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.22'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.79%  '
$ws.Range('E38').Value = '  +1.68%  '
$ws.Range('E39').Value = '  -4.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.105'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +20.31%  '
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '181.20'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.347'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.63%  '
$ws.Range('E45').Value = '  -1.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.19'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('E47').Value = '  -2.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.33'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.573'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.660'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.12%  '
$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.72'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.03%  '
